$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# name column (A) for rows 2 and 3
$ws.Range("A2").Value = "Harry Ried"
$ws.Range("A3").Value = "Adolf Mueller"

# pan column (B) for rows 2 and 3
$ws.Range("B2").Value = "CDBC4565AS"
$ws.Range("B3").Value = "ERTY7831WD"

# new row 4 - name and pan
$ws.Range("A4").Value = "Karen Reed"
$ws.Range("B4").Value = "VBDG8932JK"

# street_num column (C)
$ws.Range("C2").Value = 12
$ws.Range("C3").Value = 24
$ws.Range("C4").Value = 78

# street_name column (D)
$ws.Range("D2").Value = "74th Street"
$ws.Range("D3").Value = "4th Street"

# house_num column (E)
$ws.Range("E2").Value = 61
$ws.Range("E3").Value = ""

# locality column (F)
$ws.Range("F2").Value = "Peking nagar"
$ws.Range("F3").Value = "Urban county"

# city column (G)
$ws.Range("G2").Value = "Moonbase"
$ws.Range("G3").Value = "Marsbase"

# state column (H)
$ws.Range("H2").Value = "Ontairo"
$ws.Range("H3").Value = "Kemp"

# pin column (I)
$ws.Range("I2").Value = 456123
$ws.Range("I3").Value = 487612

# portalpass column (J) - reuse existing shared strings
$ws.Range("J2").Value = "abcdrt"
$ws.Range("J3").Value = "djfhtedpdk"

# row 4 remaining columns
$ws.Range("D4").Value = "56th Street"
$ws.Range("E4").Value = 21
$ws.Range("F4").Value = "Rural county"
$ws.Range("G4").Value = "Ergocenter"
$ws.Range("H4").Value = "Kemp"
$ws.Range("I4").Value = 125753
$ws.Range("J4").Value = "gfkhasgka"

$ws.Range("C5").Select()
